# Fill in the previously-blank row 5 ("pivot" #5) joint-data entry on the
# Inputs sheet, and move the active selection to B8 (matches the updated
# sheetView selection in the saved file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")
$ws.Activate()

$ws.Range("B7").Value = "r"
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 25
$ws.Range("F7").Value = "coupler"

$ws.Range("B8").Select()
